$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in C1 from "Tenant_ID" to "Tenant Passport ID Number"
$ws.Range("C1").Value = "Tenant Passport ID Number"

# Update the active selection to C2 (was D4)
$ws.Range("C2").Select()
